# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. right before
#    the existing "2022-Q3" sheet), and fill it with the quarterly fund
#    holdings table (same shape as the other quarter sheets).
# 2. Update the "总计" (summary) sheet: insert a new row right under the
#    header for 2022-Q4, push the previously existing quarters down by one
#    row, and keep the running index column (A) sequential (0..7).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet, positioned before "2022-Q3"
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Header row (B1:H1) - bold, centered, top-aligned, thin border
$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $q4Headers) {
    $q4.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Data rows: index, fund code, fund name, fund size, total stock position,
# position ratio, holding market value (100M yuan), position rank
$q4Data = @(
    @(0, "006218", "富国生物医药科技混合A", "8.81", "91.04", "4.87", "0.4290", 6),
    @(1, "100016", "富国天源沪港深平衡混合A", "5.01", "67.50", "4.06", "0.2034", 4),
    @(2, "011308", "富国生物医药科技混合C", "1.93", "91.04", "4.87", "0.0940", 6),
    @(3, "015228", "华夏创新研选混合C", "1.08", "92.69", "3.47", "0.0375", 10),
    @(4, "519097", "新华中小市值优选混合", "0.66", "70.51", "2.49", "0.0164", 9),
    @(5, "015227", "华夏创新研选混合A", "0.46", "92.69", "3.47", "0.0160", 10),
    @(6, "014931", "富国天源沪港深平衡混合C", "0.03", "67.50", "4.06", "0.0012", 4)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]

    $codeCell = $q4.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    $q4.Cells.Item($r, 3).Value = $row[2]

    $sizeCell = $q4.Cells.Item($r, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $row[3]

    $posCell = $q4.Cells.Item($r, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $row[4]

    $ratioCell = $q4.Cells.Item($r, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $row[5]

    $mvCell = $q4.Cells.Item($r, 7)
    $mvCell.NumberFormat = "@"
    $mvCell.Value = $row[6]

    $q4.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# Apply the bold / centered / top-aligned / thin-bordered look used by the
# other quarter sheets. Style the first header cell directly (this is the
# "seed" that defines the look), then stamp every other header cell and
# every index-column cell (column A, rows 2-8) with the exact same format
# by copy/pasting formats from that seed cell within this sheet - this
# keeps every styled cell on this new sheet mapped to one single style.
$seed = $q4.Cells.Item(1, 2)
$seed.Font.Bold = $true
$seed.HorizontalAlignment = -4108
$seed.VerticalAlignment = -4160
$seed.Borders.LineStyle = 1
$seed.Copy()

$q4.Range("C1:H1").PasteSpecial(-4122)
$q4.Range("A2:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q4 row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Existing quarters (previously rows 2-8, now shifted to rows 3-9),
# re-read from the sheet itself before we overwrite anything so we don't
# have to retype the historical numbers. (.Value2 - not .Value - is used
# for reading: the COM shim's .Value getter does not resolve correctly.)
$existing = @()
for ($row = 2; $row -le 8; $row++) {
    $label = $total.Cells.Item($row, 2).Value2
    $count = $total.Cells.Item($row, 3).Value2
    $mv = $total.Cells.Item($row, 4).Value2
    $existing += , @($label, $count, $mv)
}

# A cell that already carries the bold/centered/bordered "index column"
# look, used below as the format source for every other cell in column A.
$aStyleSrc = $total.Cells.Item(2, 1)

# Shift the existing quarters down by one row (row 8 -> 9, ..., row 2 -> 3)
for ($i = $existing.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $vals = $existing[$i]

    $idxCell = $total.Cells.Item($destRow, 1)
    $idxCell.Value = $i + 1
    $aStyleSrc.Copy()
    $idxCell.PasteSpecial(-4122)

    $lblCell = $total.Cells.Item($destRow, 2)
    $lblCell.ClearFormats()
    $lblCell.Value = $vals[0]

    $cntCell = $total.Cells.Item($destRow, 3)
    $cntCell.ClearFormats()
    $cntCell.Value = $vals[1]

    $mvCell2 = $total.Cells.Item($destRow, 4)
    $mvCell2.ClearFormats()
    $mvCell2.Value = $vals[2]
}
$excel.CutCopyMode = $false

# New 2022-Q4 row, inserted at row 2
$newIdxCell = $total.Cells.Item(2, 1)
$newIdxCell.Value = 0
$aStyleSrc.Copy()
$newIdxCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newLblCell = $total.Cells.Item(2, 2)
$newLblCell.ClearFormats()
$newLblCell.Value = "2022-Q4"

$newCntCell = $total.Cells.Item(2, 3)
$newCntCell.ClearFormats()
$newCntCell.Value = 7

$newMvCell = $total.Cells.Item(2, 4)
$newMvCell.ClearFormats()
$newMvCell.Value = 0.8

$total.Range("A1").Select()
$total.Activate()
